$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.237.31"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.710.34"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.87"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "660.58"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.424"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.07"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.708.81"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000321"
$ws.Range("E12").Value = "  +20.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.39"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.86"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.401.27"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.019.28"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.09"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.704.49"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.07"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.70"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.505"
$ws.Range("E22").Value = "  -3.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "520.99"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.44"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000219"
$ws.Range("E25").Value = "  +7.33%  "
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.41"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.195"
$ws.Range("E28").Value = "  +15.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.54"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.98"
$ws.Range("E30").Value = "  +4.01%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "656.10"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.27"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.594"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.86"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.86"
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("E43").Value = "  +3.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.489"
$ws.Range("E44").Value = "  +8.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.20"
$ws.Range("E45").Value = "  -5.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.973"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.62"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.75"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.49"
$ws.Range("E51").Value = "  +1.14%  "
